$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 54-73 (shrinks used range to A1:B53)
$ws.Rows("54:73").Delete()

# Update data rows 2-53 with new date serials (A) and values (B)
$ws.Cells.Item(2, 1).Value = 39400
$ws.Cells.Item(2, 2).Value = 3.582807286231798
$ws.Cells.Item(3, 1).Value = 39583
$ws.Cells.Item(3, 2).Value = -3.395106996021084
$ws.Cells.Item(4, 1).Value = 39765
$ws.Cells.Item(4, 2).Value = 3.499081613727355
$ws.Cells.Item(5, 1).Value = 39948
$ws.Cells.Item(5, 2).Value = 2.643706631981502
$ws.Cells.Item(6, 1).Value = 40130
$ws.Cells.Item(6, 2).Value = 2.323758253012315
$ws.Cells.Item(7, 1).Value = 40310
$ws.Cells.Item(7, 2).Value = 1.86889045152472
$ws.Cells.Item(8, 1).Value = 40494
$ws.Cells.Item(8, 2).Value = 1.711273044056469
$ws.Cells.Item(9, 1).Value = 40676
$ws.Cells.Item(9, 2).Value = 1.792365984063807
$ws.Cells.Item(10, 1).Value = 40862
$ws.Cells.Item(10, 2).Value = 1.894327589997687
$ws.Cells.Item(11, 1).Value = 41044
$ws.Cells.Item(11, 2).Value = -1.203281951233052
$ws.Cells.Item(12, 1).Value = 41228
$ws.Cells.Item(12, 2).Value = 1.349970205759888
$ws.Cells.Item(13, 1).Value = 41409
$ws.Cells.Item(13, 2).Value = -1.129543206172372
$ws.Cells.Item(14, 1).Value = 41592
$ws.Cells.Item(14, 2).Value = 1.457269272787158
$ws.Cells.Item(15, 1).Value = 41774
$ws.Cells.Item(15, 2).Value = -0.745755750114057
$ws.Cells.Item(16, 1).Value = 41957
$ws.Cells.Item(16, 2).Value = 0.06404262623020429
$ws.Cells.Item(17, 1).Value = 42137
$ws.Cells.Item(17, 2).Value = 0.2463991311210521
$ws.Cells.Item(18, 1).Value = 42321
$ws.Cells.Item(18, 2).Value = 0.6724301216571575
$ws.Cells.Item(19, 1).Value = 42503
$ws.Cells.Item(19, 2).Value = 0.6559243910605232
$ws.Cells.Item(20, 1).Value = 42689
$ws.Cells.Item(20, 2).Value = 0.593607288163227
$ws.Cells.Item(21, 1).Value = 42867
$ws.Cells.Item(21, 2).Value = 0.3249989166702818
$ws.Cells.Item(22, 1).Value = 43053
$ws.Cells.Item(22, 2).Value = 0.8673551286853183
$ws.Cells.Item(23, 1).Value = 43145
$ws.Cells.Item(23, 2).Value = 0.6331942894404392
$ws.Cells.Item(24, 1).Value = 43235
$ws.Cells.Item(24, 2).Value = 1.523252944018139
$ws.Cells.Item(25, 1).Value = 43326
$ws.Cells.Item(25, 2).Value = -0.5201366209837346
$ws.Cells.Item(26, 1).Value = 43418
$ws.Cells.Item(26, 2).Value = 0.3
$ws.Cells.Item(27, 1).Value = 43510
$ws.Cells.Item(27, 2).Value = 0.1895219545376108
$ws.Cells.Item(28, 1).Value = 43600
$ws.Cells.Item(28, 2).Value = -1.1
$ws.Cells.Item(29, 1).Value = 43691
$ws.Cells.Item(29, 2).Value = -0.558617334364854
$ws.Cells.Item(30, 1).Value = 43783
$ws.Cells.Item(30, 2).Value = 0.9418773066947779
$ws.Cells.Item(31, 1).Value = 43875
$ws.Cells.Item(31, 2).Value = 1.001530035891491
$ws.Cells.Item(32, 1).Value = 43966
$ws.Cells.Item(32, 2).Value = -0.3
$ws.Cells.Item(33, 1).Value = 44068
$ws.Cells.Item(33, 2).Value = 0.5863322451846074
$ws.Cells.Item(34, 1).Value = 44159
$ws.Cells.Item(34, 2).Value = 1.879044851730669
$ws.Cells.Item(35, 1).Value = 44251
$ws.Cells.Item(35, 2).Value = 0.5891195155317774
$ws.Cells.Item(36, 1).Value = 44341
$ws.Cells.Item(36, 2).Value = -1.620076036519961
$ws.Cells.Item(37, 1).Value = 44432
$ws.Cells.Item(37, 2).Value = 1.285713798234809
$ws.Cells.Item(38, 1).Value = 44525
$ws.Cells.Item(38, 2).Value = 0.9260365910423758
$ws.Cells.Item(39, 1).Value = 44617
$ws.Cells.Item(39, 2).Value = 1.323876910632023
$ws.Cells.Item(40, 1).Value = 44706
$ws.Cells.Item(40, 2).Value = -2.136777872354216
$ws.Cells.Item(41, 1).Value = 44798
$ws.Cells.Item(41, 2).Value = 0.7934051314956321
$ws.Cells.Item(42, 1).Value = 44890
$ws.Cells.Item(42, 2).Value = 1.029975481959951
$ws.Cells.Item(43, 1).Value = 44981
$ws.Cells.Item(43, 2).Value = 0.7043810499788776
$ws.Cells.Item(44, 1).Value = 45071
$ws.Cells.Item(44, 2).Value = 0.1406555772994125
$ws.Cells.Item(45, 1).Value = 45163
$ws.Cells.Item(45, 2).Value = 0.4268111723034735
$ws.Cells.Item(46, 1).Value = 45254
$ws.Cells.Item(46, 2).Value = 0.3049656067824742
$ws.Cells.Item(47, 1).Value = 45345
$ws.Cells.Item(47, 2).Value = 0.7386451510207621
$ws.Cells.Item(48, 1).Value = 45436
$ws.Cells.Item(48, 2).Value = -1.147231087866345
$ws.Cells.Item(49, 1).Value = 45534
$ws.Cells.Item(49, 2).Value = 0.5563279776301044
$ws.Cells.Item(50, 1).Value = 45618
$ws.Cells.Item(50, 2).Value = 0.6119095091543301
$ws.Cells.Item(51, 1).Value = 45713
$ws.Cells.Item(51, 2).Value = -0.1750403534197602
$ws.Cells.Item(52, 1).Value = 45800
$ws.Cells.Item(52, 2).Value = 2.014657057377804
$ws.Cells.Item(53, 1).Value = 45891
$ws.Cells.Item(53, 2).Value = 1.176666004305858
